$wb = $excel.ActiveWorkbook

$wsScript = $wb.Worksheets.Item("script")
$wsStudy = $wb.Worksheets.Item("study")
$wsIdMgmt = $wb.Worksheets.Item("id_management")
$wsReq = $wb.Worksheets.Item("requirements")

# -----------------------------------------------------------------
# 1) "script" sheet: add new script row (row 5) -- new R script for
#    including MCO stays by diagnosis.
# -----------------------------------------------------------------
$wsScript.Cells.Item(5, 3).Value = "inclusion_sej_mco_par_diag.R"
$wsScript.Cells.Item(5, 4).Value = "../scripts/R/inclusion_sej_mco_par_diag.R"
$wsScript.Cells.Item(5, 8).Value = "Inclure tous les séjours en MCO avec un ou plusieurs diagnostics sur plusieurs années"
$wsScript.Cells.Item(5, 2).Value = "Inclusion de séjours en MCO par les diagnostics"

$wsScript.Cells.Item(5, 1).Value = 4
$wsScript.Cells.Item(5, 5).Value = 1
$wsScript.Cells.Item(5, 6).Value = 1
$wsScript.Cells.Item(5, 9).Value = "R"

# date (publication_date) - copy number format from the cell above so it
# reuses the existing date style instead of minting a new number format
$wsScript.Cells.Item(4, 7).Copy()
$wsScript.Cells.Item(5, 7).PasteSpecial(-4122)
$wsScript.Cells.Item(5, 7).Value = 45828

# script_path column on row 5 carries the same font override style as row 4
$wsScript.Cells.Item(4, 4).Copy()
$wsScript.Cells.Item(5, 4).PasteSpecial(-4122)
$wsScript.Cells.Item(5, 4).Value = "../scripts/R/inclusion_sej_mco_par_diag.R"

# column widths (script_name, script_path, script_description got widened)
$wsScript.Columns.Item(3).ColumnWidth = 29.53
$wsScript.Columns.Item(4).ColumnWidth = 26.17
$wsScript.Columns.Item(8).ColumnWidth = 18.85

# -----------------------------------------------------------------
# 2) "study" sheet: rewrite study #1 + add new study #3 (catatonia)
# -----------------------------------------------------------------
$wsStudy.Cells.Item(2, 3).Value = "Aurélie Lescroart"
$wsStudy.Cells.Item(2, 2).Value = "Respect des recommandations de primo-prescription des antidépresseurs en France : étude rétrospective observationnelle à partir du SNDS entre 2012 et 2022 "
$wsStudy.Cells.Item(2, 7).Value = "L’étude décrit les séquences de traitement par antidépresseurs avec un suivi de 2 ans dans le cadre d’un premier épisode dépressif entre 2012 et 2022, ainsi que leur conformité aux recommandations de la Haute Autorité de Santé. Une comparaison est réalisée entre les patients suivis uniquement par leur médecin généraliste et ceux ayant reçu au moins une prescription d’antidépresseur de la part d’un psychiatre au cours de leur séquence de traitement."
$wsStudy.Cells.Item(2, 5).Value = 2
$wsStudy.Cells.Item(2, 6).Value = 5

$wsStudy.Cells.Item(4, 2).Value = "Mortalité chez les patients catatoniques"
$wsStudy.Cells.Item(4, 3).Value = "Ali Amad"
$wsStudy.Cells.Item(4, 7).Value = "L'objectif principale de l'étude est de décrire la mortalité chez les patients catatoniques avec une analyse de survie. L'objectif secondaire est d'évaluer l'impact de l'ECT sur ette mortalité."
$wsStudy.Cells.Item(4, 1).Value = 3
$wsStudy.Cells.Item(4, 5).Value = 5
$wsStudy.Cells.Item(4, 6).Value = 1

# publication date on the new study row reuses the date style from D3
$wsStudy.Cells.Item(3, 4).Copy()
$wsStudy.Cells.Item(4, 4).PasteSpecial(-4122)
$wsStudy.Cells.Item(4, 4).Value = 45828

# -----------------------------------------------------------------
# 3) "id_management" sheet: script <-> study links, updated for the
#    new script/study and a corrected mapping; one row net removed.
# -----------------------------------------------------------------
$wsIdMgmt.Cells.Item(2, 1).Value = 2
$wsIdMgmt.Cells.Item(2, 2).Value = 1
$wsIdMgmt.Cells.Item(3, 1).Value = 3
$wsIdMgmt.Cells.Item(3, 2).Value = 1
$wsIdMgmt.Cells.Item(4, 1).Value = 2
$wsIdMgmt.Cells.Item(4, 2).Value = 2
$wsIdMgmt.Cells.Item(5, 1).Value = 4
$wsIdMgmt.Cells.Item(5, 2).Value = 3
$wsIdMgmt.Rows.Item(6).Delete()

# -----------------------------------------------------------------
# 4) "requirements" sheet: new script (#4) requires script #1
# -----------------------------------------------------------------
$wsReq.Cells.Item(4, 1).Value = 4
$wsReq.Cells.Item(4, 2).Value = 1

# -----------------------------------------------------------------
# 5) leftover selection / active-sheet UI state from the edit session
# -----------------------------------------------------------------
$wsScript.Range("A6").Select()
$wsStudy.Range("A5").Select()
$wsReq.Range("B4").Select()
$wsIdMgmt.Activate()
$wsIdMgmt.Range("D14").Select()
